## ------------------------------------------------------------------------
## penalties_data.xlsx edit
##  1. Rename "late filing" -> "£100 late filing" and tweak its decile
##     ranking numbers (B7:B12, E12).
##  2. Duplicate that sheet to build a brand-new "£300 late filing" sheet
##     (same decile/rank columns, new £ amounts in the value columns,
##     trailing L:M group dropped).
##  3. Apply the same decile-ranking tweak to "late payment" and update its
##     selection.
##  4. "personal allowance" is untouched.
## ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "late filing" -> "£100 late filing"
# ---------------------------------------------------------------------
$ws100 = $wb.Worksheets.Item(1)
$ws100.Name = "£100 late filing"

$ws100.Range("B7").Value  = 17
$ws100.Range("B8").Value  = 22
$ws100.Range("B9").Value  = 29
$ws100.Range("B10").Value = 39
$ws100.Range("B11").Value = 51
$ws100.Range("B12").Value = 84
$ws100.Range("E12").Value = 87

# ---------------------------------------------------------------------
# 2. New "£300 late filing" sheet, built from a copy of "£100 late filing"
# ---------------------------------------------------------------------
$ws100.Copy($null, $ws100)
$ws300 = $wb.Worksheets.Item(2)
$ws300.Name = "£300 late filing"

# drop the 4th (2021/22) decile-value group entirely
$ws300.Range("L1:M12").EntireColumn.Delete()
# keep column K (spacer) but blank out its contents/header
$ws300.Range("K1:K12").ClearContents()

# 2018/19 group
$c = New-Object 'object[,]' 10,2
$c[0,0]=29000 ; $c[0,1]=13000
$c[1,0]=11000 ; $c[1,1]=4000
$c[2,0]=8000  ; $c[2,1]=3000
$c[3,0]=7000  ; $c[3,1]=3000
$c[4,0]=6000  ; $c[4,1]=2000
$c[5,0]=6000  ; $c[5,1]=3000
$c[6,0]=5000  ; $c[6,1]=2000
$c[7,0]=5000  ; $c[7,1]=2000
$c[8,0]=7000  ; $c[8,1]=2000
$c[9,0]=9000  ; $c[9,1]=3000
$ws300.Range("C3:D12").Value = $c

# 2019/20 group
$f = New-Object 'object[,]' 10,2
$f[0,0]=31000 ; $f[0,1]=14000
$f[1,0]=12000 ; $f[1,1]=5000
$f[2,0]=9000  ; $f[2,1]=3000
$f[3,0]=8000  ; $f[3,1]=3000
$f[4,0]=8000  ; $f[4,1]=3000
$f[5,0]=7000  ; $f[5,1]=3000
$f[6,0]=7000  ; $f[6,1]=3000
$f[7,0]=6000  ; $f[7,1]=2000
$f[8,0]=8000  ; $f[8,1]=3000
$f[9,0]=11000 ; $f[9,1]=4000
$ws300.Range("F3:G12").Value = $f

# 2020/21 group
$i = New-Object 'object[,]' 10,2
$i[0,0]=17000 ; $i[0,1]=7000
$i[1,0]=8000  ; $i[1,1]=3000
$i[2,0]=6000  ; $i[2,1]=2000
$i[3,0]=5000  ; $i[3,1]=1000
$i[4,0]=5000  ; $i[4,1]=2000
$i[5,0]=4000  ; $i[5,1]=1000
$i[6,0]=4000  ; $i[6,1]=1000
$i[7,0]=3000  ; $i[7,1]=1000
$i[8,0]=4000  ; $i[8,1]=2000
$i[9,0]=5000  ; $i[9,1]=1000
$ws300.Range("I3:J12").Value = $i

# comma-style, no-decimal number format on the 2018/19 amount columns
$ws300.Range("C3:D12").Style = "Comma"
$ws300.Range("C3:D12").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

$ws300.Range("E17").Select()

# ---------------------------------------------------------------------
# 3. "late payment" decile-ranking tweak + selection
# ---------------------------------------------------------------------
$wsPayment = $wb.Worksheets.Item("late payment")

$wsPayment.Range("B7").Value  = 17
$wsPayment.Range("B8").Value  = 22
$wsPayment.Range("B9").Value  = 29
$wsPayment.Range("B10").Value = 39
$wsPayment.Range("B11").Value = 51
$wsPayment.Range("B12").Value = 84
$wsPayment.Range("E12").Value = 87

$wsPayment.Range("E3:E12").Select()

# ---------------------------------------------------------------------
# Leave the £300 sheet as the active tab, matching the source edit
# ---------------------------------------------------------------------
$ws300.Activate()
$ws300.Range("E17").Select()
